# Release v0.6.0 update.
#
# 1) A4 / A5 had a redundant "apply fill" style left over from earlier
#    editing; clear it back to the default (General/no style) look.
# 2) Rows 17 and 18 had their transaction data entered in the wrong order
#    (a Poloniex->Desktop-wallet withdrawal/deposit pair) - swap the two
#    rows' data back (everything except the Timestamp column L, which
#    stays put on its own row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Clear the stray fill-applying style on A4/A5 -----------------
$ws.Range("A4").ClearFormats()
$ws.Range("A5").ClearFormats()

# --- 2) Swap rows 17 & 18 (all columns except L, the Timestamp) ------
$ws.Range("A17").Value = "Withdrawal"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = ""
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = "BTC"
$ws.Range("H17").Value = 0.0001259
$ws.Range("I17").Value = "BTC"
$ws.Range("K17").Value = "Desktop wallet"
$ws.Range("M17").Value = "from Poloniex"

$ws.Range("A18").Value = "Deposit"
$ws.Range("B18").Value = 6
$ws.Range("C18").Value = "BTC"
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""
$ws.Range("K18").Value = "Poloniex"
$ws.Range("M18").Value = "to Desktop wallet"

# Leave the selection where the author last left it (P23).
$ws.Range("P23").Select() | Out-Null
